$wb = $excel.ActiveWorkbook

# The workbook currently ends with sheet "07.09" (the most recent daily
# ranking tab). This release adds the next day's ranking tab ("08.09"),
# built the same way all the other daily tabs were: duplicate the most
# recent tab, rename it, and update its title cell / sort record.

$sheetCount = $wb.Worksheets.Count
$sourceSheet = $wb.Worksheets.Item($sheetCount)

# Duplicate "07.09" and place the copy right after it -> becomes the new
# last sheet and, as in Excel, the newly created/copied sheet becomes the
# active tab while the previous last sheet loses its "selected" state.
$sourceSheet.Copy($null, $sourceSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "08.09"

# Update the title cell (A1) to reference the new date.
$newSheet.Range("A1").Value = "Рейтинг подразделений по 4 показателям за 08.09.2022"

# Re-run the same ranking sort (ascending by rank, column A) on the new
# sheet so it carries its own freshly generated sort record instead of
# simply inheriting the one copied from "07.09".
$sortSpec = $newSheet.Sort
$sortSpec.SortFields.Clear()
$sortSpec.SortFields.Add($newSheet.Range("A3"))
$sortSpec.SetRange($newSheet.Range("A3:B15"))
$sortSpec.Header = 0
$sortSpec.Apply()
